$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 364.1
$ws.Range("I12").Value = 313
$ws.Range("J12").Value = 415.2
$ws.Range("K12").Value = 313
$ws.Range("L12").Value = 415.2
$ws.Range("M12").Value = -143
$ws.Range("N12").Value = -755.2

$ws.Range("H33").Value = 90910090
$ws.Range("I33").Value = 128
$ws.Range("J33").Value = 200002030
$ws.Range("K33").Value = 128
$ws.Range("L33").Value = 200002030
$ws.Range("M33").Value = 101
$ws.Range("N33").Value = -200002488

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H76").Value = 20666.334
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 20666.334
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 20666.334
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -21296.334

$ws.Range("H79").Value = 20666.334
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 20666.334
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 20666.334
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -22850.334

$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2430.5938
$ws.Range("I32").Value = 945.71155
$ws.Range("J32").Value = 8865.083000000001
$ws.Range("K32").Value = 945.71155
$ws.Range("L32").Value = 8865.083000000001
$ws.Range("M32").Value = -658.71155
$ws.Range("N32").Value = -9439.083000000001

$ws.Range("H74").Value = 1666.6538
$ws.Range("I74").Value = 1514.3684
$ws.Range("J74").Value = 2080
$ws.Range("K74").Value = 1514.3684
$ws.Range("L74").Value = 2080
$ws.Range("M74").Value = -640.3684000000001
$ws.Range("N74").Value = -3828

$ws.Range("H77").Value = 1666.6538
$ws.Range("I77").Value = 1514.3684
$ws.Range("J77").Value = 2080
$ws.Range("K77").Value = 7571.842000000001
$ws.Range("L77").Value = 10400
$ws.Range("M77").Value = -3203.842000000001
$ws.Range("N77").Value = -19136

$ws.Range("H102").Value = 66429.5
$ws.Range("I102").Value = 22384.666
$ws.Range("K102").Value = 22384.666
$ws.Range("M102").Value = -20762.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2814.35
$ws.Range("I20").Value = 2455.75
$ws.Range("K20").Value = 2455.75
$ws.Range("M20").Value = -2208.75

$ws.Range("H86").Value = 2749
$ws.Range("I86").Value = 2446.6
$ws.Range("K86").Value = 2446.6
$ws.Range("M86").Value = -1323.6

$ws.Range("H89").Value = 2749
$ws.Range("I89").Value = 2446.6
$ws.Range("K89").Value = 12233
$ws.Range("M89").Value = -6617

$ws.Range("H134").Value = 4686.9287
$ws.Range("I134").Value = 4326.5835
$ws.Range("J134").Value = 6849
$ws.Range("K134").Value = 12979.7505
$ws.Range("L134").Value = 20547
$ws.Range("M134").Value = -10444.7505
$ws.Range("N134").Value = -25617

$ws.Range("H141").Value = 244999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 244999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 244999.5
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -255359.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4001970
$ws.Range("I19").Value = 5002250
$ws.Range("K19").Value = 5002250
$ws.Range("M19").Value = -5002080

$ws.Range("H24").Value = 4001970
$ws.Range("I24").Value = 5002250
$ws.Range("K24").Value = 5002250
$ws.Range("M24").Value = -5002080

$ws.Range("H31").Value = 2053.8215
$ws.Range("I31").Value = 2524.75
$ws.Range("K31").Value = 2524.75
$ws.Range("M31").Value = -2229.75

$ws.Range("H34").Value = 2053.8215
$ws.Range("I34").Value = 2524.75
$ws.Range("K34").Value = 2524.75
$ws.Range("M34").Value = -2322.75

$ws.Range("H134").Value = 3001.9092
$ws.Range("I134").Value = 2918.1667
$ws.Range("J134").Value = 3102.4
$ws.Range("K134").Value = 8754.500100000001
$ws.Range("L134").Value = 9307.200000000001
$ws.Range("M134").Value = -6219.500100000001
$ws.Range("N134").Value = -14377.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1703.6875
$ws.Range("I122").Value = 2294.5
$ws.Range("K122").Value = 20650.5
$ws.Range("M122").Value = -18200.5

$ws.Range("H129").Value = 36113224
$ws.Range("J129").Value = 20004000
$ws.Range("L129").Value = 60012000
$ws.Range("N129").Value = -60022000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5618.278
$ws.Range("I80").Value = 5365.8887
$ws.Range("J80").Value = 5870.6665
$ws.Range("K80").Value = 5365.8887
$ws.Range("L80").Value = 5870.6665
$ws.Range("M80").Value = -4367.8887
$ws.Range("N80").Value = -7866.6665

$ws.Range("H83").Value = 5618.278
$ws.Range("I83").Value = 5365.8887
$ws.Range("J83").Value = 5870.6665
$ws.Range("K83").Value = 26829.4435
$ws.Range("L83").Value = 29353.3325
$ws.Range("M83").Value = -21837.4435
$ws.Range("N83").Value = -39337.3325

$ws.Range("H102").Value = 3742.2273
$ws.Range("I102").Value = 3616.5
$ws.Range("K102").Value = 3616.5
$ws.Range("M102").Value = -1994.5

$ws.Range("H113").Value = 2191
$ws.Range("J113").Value = 2191
$ws.Range("L113").Value = 2191
$ws.Range("N113").Value = -6531

$ws.Range("H132").Value = 3924
$ws.Range("I132").Value = 4075.6
$ws.Range("J132").Value = 3671.3333
$ws.Range("K132").Value = 12226.8
$ws.Range("L132").Value = 11013.9999
$ws.Range("M132").Value = -9696.799999999999
$ws.Range("N132").Value = -16073.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 751.2941
$ws.Range("I16").Value = 720.1875
$ws.Range("J16").Value = 1249
$ws.Range("K16").Value = 720.1875
$ws.Range("L16").Value = 1249
$ws.Range("M16").Value = -550.1875
$ws.Range("N16").Value = -1589

$ws.Range("H61").Value = 13999.667
$ws.Range("I61").Value = 19999.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 19999.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -19797.5
$ws.Range("N61").Value = -2404

$ws.Range("H68").Value = 1903.3125
$ws.Range("I68").Value = 1696.8667
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 1696.8667
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -947.8667
$ws.Range("N68").Value = -6498

$ws.Range("H71").Value = 1903.3125
$ws.Range("I71").Value = 1696.8667
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 8484.333500000001
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -4740.333500000001
$ws.Range("N71").Value = -32488

$ws.Range("H113").Value = 13999.667
$ws.Range("I113").Value = 19999.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 19999.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -17829.5
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3845.3333
$ws.Range("I81").Value = 1324.8334
$ws.Range("K81").Value = 2649.6668
$ws.Range("M81").Value = -1588.6668

$ws.Range("H84").Value = 3845.3333
$ws.Range("I84").Value = 1324.8334
$ws.Range("K84").Value = 13248.334
$ws.Range("M84").Value = -7944.333999999999

$ws.Range("H123").Value = 42849.6
$ws.Range("I123").Value = 19999.5
$ws.Range("K123").Value = 19999.5
$ws.Range("M123").Value = -15099.5

$ws.Range("H136").Value = 3882.8333
$ws.Range("I136").Value = 3882.8333
$ws.Range("K136").Value = 11648.4999
$ws.Range("M136").Value = -9098.499899999999

